# "change tracing strategy and save wallet labels"
#
# Append the newest Date/USDValue observation as row 9:
#   A9 = "2024-10-05"  (stored as text, like the other Date cells)
#   B9 = 6.209E-05

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates formatted as plain text (e.g. "2024-08-25"), not
# real Excel dates. Assigning the literal string straight to A9.Value
# would make Excel's "looks like a date" literal-entry heuristic kick in
# and silently store a date serial number instead of the text. To avoid
# that (and to avoid leaving unused NumberFormat/style definitions behind
# from a NumberFormat="@" workaround), build the text via a formula in a
# scratch cell, copy it, and paste only the resulting value into A9 -
# pasted values are never reinterpreted as dates.
$ws.Range("D1").Formula = "=""2024-10-05"""
$ws.Range("D1").Copy()
$ws.Range("A9").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("D1").Clear()

$ws.Range("B9").Value = 0.00006209
